$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Average of column J (the ROI k value) across rows 2..11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$rng = $ws.Range("B14:B17")
$rng.Font.Bold = $true
$rng.Font.Size = 12
$rng.VerticalAlignment = -4108

$ws.Range("A14:B17").RowHeight = 15.6

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select()
